# "I swapped 4o to 4.1mini in results"
#
# The ResultFigures sheet lists model results in 4 groups of 3 rows each
# (GPT-5-mini / GPT-4o-mini / o3-mini) under columns J (model name) and
# K (percent value), for question types Recall, Modus Ponens, Modus
# Tollens and ... (rows 6-8, 12-14, 18-20, 24-26). Rename every
# "GPT-4o-mini" label to "GPT-4.1-mini" and update the three K values
# that changed for that model's row in each group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 group (Recall): label rename only, value unchanged (20).
$ws.Range("J7").Value = "GPT-4.1-mini"

# Row 13 group (Modus Ponens): label rename + value 40 -> 80.
$ws.Range("J13").Value = "GPT-4.1-mini"
$ws.Range("K13").Value = 80

# Row 19 group (Modus Tollens): label rename + value 60 -> 80.
$ws.Range("J19").Value = "GPT-4.1-mini"
$ws.Range("K19").Value = 80

# Row 25 group (last question type): label rename + value 80 -> 60.
$ws.Range("J25").Value = "GPT-4.1-mini"
$ws.Range("K25").Value = 60

# Match the saved selection state left behind by the edit.
$ws.Range("K14").Select()
